$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a whole new column before column B. This shifts column widths
#    and all cell data in columns B..L to C..M (for every row), and
#    auto-extends the A2:L2 title merge to A2:M2.
$ws.Columns("B:B").Insert()

# 2. Fix up the newly inserted column B for the label rows (4-8): the
#    insert copies formatting from column A (bold box) into the new
#    column B, but it must look like the rest of the input box (the
#    plain bordered style already used by C5:E5 etc). Copy that format
#    from column C onto column B for each of those rows.
$ws.Range("C4").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(4122)
$ws.Range("C5").Copy() | Out-Null
$ws.Range("B5").PasteSpecial(4122)
$ws.Range("C6").Copy() | Out-Null
$ws.Range("B6").PasteSpecial(4122)
$ws.Range("C7").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(4122)
$ws.Range("C8").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(4122)
$ws.Application.CutCopyMode = $false

# 3. New "SUCURSAL" header column in the table (write this BEFORE
#    "EMPRESA :" below so the shared-string table picks up the same
#    ordering as the authored workbook).
$ws.Range("B10").Value = "SUCURSAL"

# 4. Row 4 becomes a two-field row: "EMPRESA :" (existing label spot)
#    plus a second "ESTABLECIMEINTO :" field built from new cells
#    G4:J4 (which inherit the look of the other input boxes).
$ws.Range("A4").Value = "EMPRESA :"

$ws.Range("G4").Value = $ws.Range("C4").Value
$ws.Range("G4").Copy() | Out-Null
$ws.Range("H4").PasteSpecial(4122)
$ws.Range("I4").PasteSpecial(4122)
$ws.Range("J4").PasteSpecial(4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A4").Copy() | Out-Null
$ws.Range("F4").PasteSpecial(4122)
$ws.Application.CutCopyMode = $false
$ws.Range("F4").Value = "ESTABLECIMEINTO :"

# 5. Re-merge the label input boxes (now one column wider: B:E instead
#    of the old B:D), and merge the new second field's box G4:J4.
$ws.Range("B4:E4").Merge()
$ws.Range("B5:E5").Merge()
$ws.Range("B6:E6").Merge()
$ws.Range("B7:E7").Merge()
$ws.Range("B8:E8").Merge()
$ws.Range("G4:J4").Merge()

# 6. Selection, matching the saved workbook view.
$ws.Range("A5").Select()
